$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 36) to the table, matching the existing rows'
# pattern: date text in column A, weekday text in column B, and two numbers.
# The date looks like a date literal, so Excel would normally auto-convert it
# to a date serial on assignment; use the leading-apostrophe text-prefix
# trick (same as typing '2025/09/29 into the cell) and then clear the
# resulting quote-prefix formatting so the cell ends up as a plain,
# unstyled text value - consistent with the rest of the column.
$ws.Range("A36").Value = "'2025/09/29"
$ws.Range("A36").ClearFormats()

$ws.Range("B36").Value = "月"
$ws.Range("C36").Value = 20
$ws.Range("D36").Value = 3
